$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the answer for "What is Dr. Huettemans favorite color?"
$ws.Range("B6").Value = "One hundred thirty seven"

# Fix the answer for "What is the capital of Pennsylvania?"
$ws.Range("B7").Value = "Harisburgh"

# Move the active selection to B6
$ws.Range("B6").Select()
